$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update "Right" count on the Marking row
$ws.Range("B11").Value = 5

# Update "Right" count and the corr/total summary on the Total row
$ws.Range("B12").Value = 95
$ws.Range("E12").Value = "95/140"
